$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 5 values
$ws.Range("A5").Value = 10076.65
$ws.Range("B5").Value = 10193.879999999999
$ws.Range("C5").Value = 19.170000000000002
$ws.Range("D5").Value = 19.39
$ws.Range("E5").Value = $true
$ws.Range("F5").Value = 1.1499999999999999
$ws.Range("G5").Value = 42609.505856481483
$ws.Range("H5").Value = $false

# Match the date/time style used by the rest of column G (style index 1 -> numFmtId 22)
$ws.Range("G4").Copy()
$ws.Range("G5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G5").Value = 42609.505856481483
$excel.CutCopyMode = $false
